$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.906.40'

$ws.Range("E2").Value = '  -2.32%  '

$ws.Range("D3").Value = '1.901.37'

$ws.Range("E3").Value = '  -4.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.59'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4588'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = '  -2.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3810'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = '  -3.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.56'
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = '  -2.51%  '

$ws.Range("E10").Value = '  -2.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9810'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.06'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  -3.56%  '

$ws.Range("D13").Value = '1.981.70'

$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("E14").Value = '  -4.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.672'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = '  -3.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07054'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = '  -1.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '84.13'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = '  -5.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009552'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '  -4.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.74'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  -4.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.005'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").Value = '28.861.15'

$ws.Range("E22").Value = '  -2.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.338'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = '  -4.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = '  -3.34%  '

$ws.Range("D25").Value = '2.152.32'

$ws.Range("E25").Value = '  -4.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.098'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.25'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.12'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = '  -3.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.587'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  -7.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.65'
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = '  -2.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.838'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = '  -5.53%  '

$ws.Range("E32").Value = '  -2.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8627'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = '  -6.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.101'
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = '  -3.36%  '

$ws.Range("E35").Value = '  -7.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.020'
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = '  -5.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05690'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  -2.81%  '

$ws.Range("E38").Value = '  -2.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.005'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = '  +0.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02036'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = '  -4.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.478'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = '  -5.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5513'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  -4.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1756'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  -4.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.318'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  -5.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.727'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  -1.10%  '

$ws.Range("E46").Value = '  -3.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.25'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = '  -6.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.104'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  -4.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06830'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  -1.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.55'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  -2.46%  '

$ws.Range("E51").Value = '  -5.67%  '
